$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the disorder name typo: "SEED" -> "SEDD"
$ws.Range("A4").Value = "SEDD"

# Update the view state to match (scrolled so row 4 is at the top, whole row 4 selected)
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("A4:XFD4").Select()
